$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: row 11 "Marking" and row 12 "Total"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -30
$ws.Range("E12").Value = "22 / 112"
